$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, shifting existing data down
$ws.Rows.Item(1).Insert()

# Populate the new header row
$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "Colour1"
$ws.Range("C1").Value = "Colour2"
$ws.Range("D1").Value = "Colour3"
$ws.Range("E1").Value = "Colour4"

# Set column widths for A:E (target stored width 17.7109375 chars / MDW7;
# this engine quantizes ColumnWidth to 1/6-character pixel steps, so the
# closest achievable stored width is 17.6666... via an input of 16.8333...)
$ws.Range("A1:E1").EntireColumn.ColumnWidth = 16.833333333333336

# Set the active cell / selection to E1
$ws.Range("E1").Select()
